$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), copying the format of the neighboring
# header cell (G1) so it keeps the bold/bordered/centered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (2..63), matching the commit's data
$saveValues = @{
    2 = 0; 3 = 0; 4 = 1; 5 = 0; 6 = 0; 7 = 1; 8 = 0; 9 = 0; 10 = 0;
    11 = 1; 12 = 0; 13 = 0; 14 = 1; 15 = 1; 16 = 0; 17 = 0; 18 = 1; 19 = 0;
    20 = 0; 21 = 1; 22 = 0; 23 = 0; 24 = 0; 25 = 1; 26 = 0; 27 = 0; 28 = 0;
    29 = 1; 30 = 0; 31 = 0; 32 = 0; 33 = 0; 34 = 0; 35 = 1; 36 = 0; 37 = 0;
    38 = 0; 39 = 1; 40 = 0; 41 = 1; 42 = 0; 43 = 1; 44 = 0; 45 = 0; 46 = 1;
    47 = 1; 48 = 1; 49 = 0; 50 = 1; 51 = 0; 52 = 0; 53 = 0; 54 = 1; 55 = 0;
    56 = 1; 57 = 0; 58 = 0; 59 = 0; 60 = 1; 61 = 0; 62 = 0; 63 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
